$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary_counts")

# Insert a new row at row 13, shifting the existing rows 13-14 down to 14-15.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new statistic.
$ws.Range("A13").Value = "Number of events with both any university response coding and any police coding"
$ws.Range("B13").Value = 360

# Update the worksheet's used-range dimension to reflect the extra row.
$ws.Range("A1:B15").Select()
